# (update) menu pasca penindakan tampilan form create
#
# - Resize the "CATATAN" table's two grid columns (table 1).
# - Change the page size from Letter (12240x15840 twips) to
#   Folio/F4 (11907x18711 twips).

$d = $word.ActiveDocument

# --- Resize columns of the first table (the "CATATAN : ${catatan}" table) ---
$catatanTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text -like "CATATAN*") {
        $catatanTable = $candidate
        break
    }
}
if ($catatanTable -eq $null) {
    $catatanTable = $d.Tables.Item(1)
}

# Widths are expressed in points (1 pt = 20 twips):
#   1838 twips -> 1811 twips (col 1)
#   7512 twips -> 7216 twips (col 2)
$catatanTable.Columns.Item(1).Width = 1811 / 20.0
$catatanTable.Columns.Item(2).Width = 7216 / 20.0

# --- Update the page size (section page setup) ---
# 12240x15840 twips (Letter, 8.5in x 11in) -> 11907x18711 twips (F4/Folio, 21cm x 33cm)
$d.PageSetup.PageWidth = 11907 / 20.0
$d.PageSetup.PageHeight = 18711 / 20.0
